$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to Text format so that numeric-looking
# strings such as "311.74" or "1.003" are preserved as literal text, matching the
# original inline-string cell content instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.468.23"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "1.824.86"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "311.74"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "0.4251"
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").Value = "0.3633"
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").Value = "0.07231"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("D10").Value = "0.8572"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("D11").Value = "20.52"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").Value = "1.804.66"
$ws.Range("E12").Value = "  -5.19%  "

$ws.Range("D13").Value = "5.414"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "6.481"
$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").Value = "0.06966"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").Value = "80.19"
$ws.Range("E17").Value = "  +1.19%  "

$ws.Range("D18").Value = "0.000008872"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").Value = "0.9995"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").Value = "15.34"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "27.387.13"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("D22").Value = "5.135"
$ws.Range("E22").Value = "  +3.23%  "

$ws.Range("D23").Value = "10.83"
$ws.Range("E23").Value = "  +5.25%  "

$ws.Range("D24").Value = "2.063.05"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").Value = "1.984"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "154.63"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").Value = "18.76"
$ws.Range("E27").Value = "  +1.52%  "

$ws.Range("D28").Value = "5.100"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("D29").Value = "113.89"
$ws.Range("E29").Value = "  -4.66%  "

$ws.Range("D30").Value = "1.803"
$ws.Range("E30").Value = "  -3.86%  "

$ws.Range("D31").Value = "0.08839"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "2.981"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.7405"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("D34").Value = "4.508"
$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").Value = "1.119"
$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").Value = "1.092"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("D38").Value = "0.05288"
$ws.Range("E38").Value = "  -2.97%  "

$ws.Range("D39").Value = "0.01925"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").Value = "2.788"
$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("D41").Value = "0.5041"
$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").Value = "0.1640"
$ws.Range("E42").Value = "  -1.22%  "

$ws.Range("D43").Value = "6.457"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").Value = "8.269"
$ws.Range("E44").Value = "  -1.34%  "

$ws.Range("D45").Value = "10.34"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "105.05"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").Value = "0.06455"
$ws.Range("E47").Value = "  -1.47%  "

$ws.Range("D48").Value = "0.4656"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").Value = "1.607"
$ws.Range("E50").Value = "  -1.64%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "63.32"
$ws.Range("E51").Value = "  -2.08%  "

# Reset the style reference on these cells back to the default "Normal" style so
# we do not leave a stray custom cell style applied (matches original formatting,
# which had no explicit style index on these data cells).
$ws.Range("D2:E51").Style = "Normal"
